# "Add UI diagrams to Appendix 4"
#
# The deck has a single slide ("Multi word phrase" diagram). This edit:
#   1. Re-colours the word "phrase" (previously part of a single run
#      " phrase") in amber (FFC000), splitting it from the leading space
#      which keeps its original background colour.
#   2. Slides the "Word 1" rectangle to the right to make room.
#   3. Slides the "Word 3" rectangle to the left and fills it with the
#      same amber (FFC000) colour used for "phrase", tying the three
#      diagram pieces together.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 1 ("Multi word phrase") -------------------------------------
# Runs today: "Multi " / "word" / " phrase"  -> characters 1-6 / 7-10 / 11-17
# Target: split the trailing " phrase" run into " " (unchanged) and
# "phrase" (amber FFC000), leaving the first two runs untouched.
$title = $s.Shapes.Item(2)
$titleRange = $title.TextFrame.TextRange
$wordRange = $titleRange.Characters(12, 6)
$wordRange.Font.Color.RGB = 49407  # 0xC0FF -> RGB(255,192,0) = FFC000

# --- Rectangle 3 ("Word 1") -----------------------------------------------
$word1 = $s.Shapes.Item(3)
$word1.Left = 244.78901647795277   # 3108820 EMU

# --- Rectangle 9 ("Word 3") -----------------------------------------------
$word3 = $s.Shapes.Item(5)
$word3.Left = 578.4769895338582    # 7346657 EMU
$word3.Fill.Solid()
$word3.Fill.ForeColor.RGB = 49407  # FFC000
